# Slide 1: every shape was selected together and dragged down-and-right by
# a uniform offset of (+555586, +833375) EMU.
#
# Shape.Left/.Top are COM `Single` (32-bit float) properties measured in
# points (1 pt = 12700 EMU); the host truncates pt*12700 to get EMU, so a
# naive point value can land 1 EMU short after the float32 round-trip. To
# hit the exact target EMU we add a tiny (0.5 EMU, i.e. 0.5/12700 pt) nudge
# before truncation.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuPerPt = 12700
$deltaXEmu = 555586
$deltaYEmu = 833375
$nudge = 0.5 / $emuPerPt

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)

    $curXEmu = [math]::Round($shape.Left * $emuPerPt)
    $curYEmu = [math]::Round($shape.Top * $emuPerPt)

    $newXEmu = $curXEmu + $deltaXEmu
    $newYEmu = $curYEmu + $deltaYEmu

    $shape.Left = ($newXEmu / $emuPerPt) + $nudge
    $shape.Top = ($newYEmu / $emuPerPt) + $nudge
}
